# fix(voters): correct Excel template headers to match validation rules
#
# New header order/text: שם, טלפון, שם משפחה, עיר, מייל
# (Required: שם, טלפון — matches backend validation; שם משפחה/עיר/מייל optional)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "שם"
$ws.Range("B1").Value = "טלפון"
$ws.Range("C1").Value = "שם משפחה"
$ws.Range("D1").Value = "עיר"
$ws.Range("E1").Value = "מייל"

# Widen the "name" columns (A: first name, C: surname, E: email) to fit
# the new name-first column order. B and D keep their original width.
$ws.Columns.Item(1).ColumnWidth = 20
$ws.Columns.Item(3).ColumnWidth = 20
$ws.Columns.Item(5).ColumnWidth = 30
